# Time recording log.xlsx - "Nädal 7" sheet updates (Repository klasside kokkuvõte, osa 25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (item 5): worked 12.03.2020, 08:00-09:00 -> 60 minutes on "Praktikum", marked done with "x"
$ws.Range("B11").Value = 43902
$ws.Range("C11").Value = 0.33333333333333331
$ws.Range("D11").Value = 0.375
$ws.Range("F11").Value = 60
$ws.Range("G11").Value = "Praktikum"
$ws.Range("I11").Value = "x"

# Row 12 (item 6): worked 12.03.2020, 09:30-09:48 -> 18 minutes on "Kodutöö 6", comment "p.25 tehtud", marked done with "x"
$ws.Range("B12").Value = 43902
$ws.Range("C12").Value = 0.39583333333333331
$ws.Range("D12").Value = 0.40833333333333338
$ws.Range("F12").Value = 18
$ws.Range("G12").Value = "Kodutöö 6"
$ws.Range("H12").Value = "p.25 tehtud"
$ws.Range("J12").Value = "x"

# Update the active selection to reflect where the author left off
$ws.Range("G15").Select()
